$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.847326993942261
$ws.Range("B1").Value = 3.084224939346313
$ws.Range("C1").Value = 2.711798667907715
$ws.Range("D1").Value = 2.948814392089844
$ws.Range("E1").Value = 2.697418928146362
